# Swap the contents of rows 3 and 4 (the "Test" and "S1" signal configs)
# so that "S1" appears in row 3 and "Test" appears in row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")

foreach ($col in $columns) {
    $val3 = $ws.Range("$col`3").Value2
    $val4 = $ws.Range("$col`4").Value2
    $ws.Range("$col`3").Value = $val4
    $ws.Range("$col`4").Value = $val3
}
